# Adds the two new "LoginManager" test-case worksheets described in the
# commit "Added unit testing for LoginManager ProcessRegistration":
#   - "LoginManager ProcessRegistratio"  (31-char truncated name) - full data
#   - "LoginManager FailedRegistration"  - header-only stub, newly started

$wb = $excel.ActiveWorkbook

# The "LoginManager ResetPassword" sheet already has the exact header layout
# we need (Test Case / Category / Partition / Test Inputs(3 cols) /
# Expected Output / Comments), so we reuse it as a formatting template.
$template = $wb.Worksheets.Item("LoginManager ResetPassword")
$lastExisting = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---------------------------------------------------------------------
# 1) "LoginManager ProcessRegistratio" - new sheet with full test data
# ---------------------------------------------------------------------
$proc = $wb.Worksheets.Add($null, $lastExisting)
$proc.Name = "LoginManager ProcessRegistratio"

# Bring over the two-row merged header (values + formatting + merged cells)
$template.Range("A1:H2").Copy($proc.Range("A1"))

# Column B ("Category") is centre-aligned in every test sheet.
$proc.Range("B3:B9").HorizontalAlignment = -4108

$procData = @(
    @(1, "valid",   "all valid inputs",    "userTwo", "Password2!",  "GetIPAddress()", $true),
    @(2, "valid",   "invalid uid format",   "user2",   "Password2!",  "GetIPAddress()", $false),
    @(3, "valid",   "invalid pwd format",   "userTwo", "passwordTwo", "GetIPAddress()", $false),
    @(4, "valid",   "empty uid",            '""',      "Password2!",  "GetIPAddress()", $false),
    @(5, "valid",   "empty pwd",            "userTwo", '""',          "GetIPAddress()", $false),
    @(6, "invalid", "null uid",             "null",    "Password2!",  "GetIPAddress()", $false),
    @(7, "invalid", "null pwd",             "userTwo", "null",        "GetIPAddress()", $false)
)

$r = 3
foreach ($row in $procData) {
    $proc.Cells.Item($r, 1).Value = $row[0]
    $proc.Cells.Item($r, 2).Value = $row[1]
    $proc.Cells.Item($r, 3).Value = $row[2]
    $proc.Cells.Item($r, 4).Value = $row[3]
    $proc.Cells.Item($r, 5).Value = $row[4]
    $proc.Cells.Item($r, 6).Value = $row[5]
    $proc.Cells.Item($r, 7).Value = $row[6]
    $r++
}

# ---------------------------------------------------------------------
# 2) "LoginManager FailedRegistration" - new sheet, header only (stub)
# ---------------------------------------------------------------------
$failed = $wb.Worksheets.Add($null, $proc)
$failed.Name = "LoginManager FailedRegistration"

$template.Range("A1:H2").Copy($failed.Range("A1"))
$failed.Range("B3:B9").HorizontalAlignment = -4108

# Leave the new sheet active, positioned where data entry would continue.
$failed.Activate()
$failed.Range("D3").Select()

$excel.CutCopyMode = $false
